$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# like "4.401" are not auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.072.79"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.661.03"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "208.18"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "0.5174"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "0.2582"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").Value = "0.06303"
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").Value = "20.94"
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "0.07537"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.664.65"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "4.401"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "0.5380"
$ws.Range("E14").Value = "  -3.77%  "
$ws.Range("D15").Value = "66.15"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "0.0₅7949"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "26.084.35"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "4.703"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "187.16"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "10.16"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "148.32"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").Value = "0.1214"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("D26").Value = "7.378"
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("D27").Value = "15.65"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "1.391"
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("D29").Value = "0.05996"
$ws.Range("E29").Value = "  -5.69%  "
$ws.Range("D30").Value = "1.262"
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("D31").Value = "3.467"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "3.398"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").Value = "1.637"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "0.9853"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").Value = "2.758"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "0.5882"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("D38").Value = "1.104.68"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "0.01594"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "5.970"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "0.8472"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "99.89"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "1.816.71"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "55.09"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("D49").Value = "0.05226"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "0.4237"
$ws.Range("E51").Value = "  -0.75%  "

Write-Host "Updated cryptos list"
